# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (want-to-go count) figures and one updated
# cover-image URL to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 5428
$ws1.Range("F6").Value  = 27
$ws1.Range("F7").Value  = 624
$ws1.Range("F8").Value  = 597
$ws1.Range("F11").Value = 1495
$ws1.Range("F12").Value = 4585
$ws1.Range("F13").Value = 442
$ws1.Range("F15").Value = 176
$ws1.Range("F17").Value = 3554
$ws1.Range("F18").Value = 180
$ws1.Range("F19").Value = 1118
$ws1.Range("F22").Value = 205
$ws1.Range("F23").Value = 26
$ws1.Range("F24").Value = 138
$ws1.Range("F25").Value = 49
$ws1.Range("F28").Value = 320
$ws1.Range("F29").Value = 34
$ws1.Range("F30").Value = 60
$ws1.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202404/N6VdMOuL1713257425864.jpeg"

# --- Sheet "全部类型" (all types, rows shifted +1 vs 展览) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value  = 5428
$ws4.Range("F7").Value  = 27
$ws4.Range("F8").Value  = 624
$ws4.Range("F9").Value  = 597
$ws4.Range("F12").Value = 1495
$ws4.Range("F13").Value = 4585
$ws4.Range("F14").Value = 442
$ws4.Range("F16").Value = 176
$ws4.Range("F18").Value = 3554
$ws4.Range("F19").Value = 180
$ws4.Range("F20").Value = 1118
$ws4.Range("F23").Value = 205
$ws4.Range("F24").Value = 26
$ws4.Range("F25").Value = 138
$ws4.Range("F26").Value = 49
$ws4.Range("F29").Value = 320
$ws4.Range("F30").Value = 34
$ws4.Range("F31").Value = 60
$ws4.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202404/N6VdMOuL1713257425864.jpeg"
